$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold numeric-looking text (e.g. "227.42"). Excel would
# silently convert a plain .Value assignment into a real number, so force
# Text format first to preserve them as strings, matching the source data.
$textCells = @("D5", "D6", "D7", "D10", "D14", "D15", "D20", "D22", "D26", "D27", "D29", "D30", "D32", "D36", "D38", "D40", "D43", "D46", "D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.719.09"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "2.026.91"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "227.42"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").Value = "0.602"
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("D7").Value = "59.81"
$ws.Range("E7").Value = "  -2.96%  "
$ws.Range("E9").Value = "  -3.68%  "
$ws.Range("D10").Value = "0.0821"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("D12").Value = "2.328.53"
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D14").Value = "21.01"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "0.768"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").Value = "1.998.87"
$ws.Range("E17").Value = "  -3.69%  "
$ws.Range("D18").Value = "37.681.28"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "5.88"
$ws.Range("E20").Value = "  -7.01%  "
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "223.55"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").Value = "167.91"
$ws.Range("D27").Value = "9.33"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("D29").Value = "18.76"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "1.25"
$ws.Range("E30").Value = "  -6.17%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "2.21"
$ws.Range("E32").Value = "  +7.61%  "
$ws.Range("E33").Value = "  -4.68%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("D36").Value = "6.41"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "3.41"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "17.96"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("D41").Value = "1.532.63"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("D43").Value = "95.45"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "1.10"
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "4.05"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "2.217.85"
$ws.Range("E51").Value = "  -1.93%  "

# Restore the default General format on those cells now that the text
# values are stored, so no stray Text-format style lingers on them.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

